$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabela1")

# Duplicate the last data row (row 16) - including its formatting - into the
# new row 17 so the added row inherits the same cell styles.
$ws.Rows("16").Copy()
$ws.Rows("17").Insert(-4121)

# Grow the table / autofilter range so row 17 becomes part of "Tabela1".
$tbl.Resize($ws.Range("A1:G17"))

# Fill in the new study-log entry for the row that was just added.
$ws.Range("A17").Value = 44844
$ws.Range("B17").Value = 0.25
$ws.Range("C17").Value = 1.2486111111111111
$ws.Range("D17").Formula = "=Tabela1[[#This Row],[HORA F]]-Tabela1[[#This Row],[HORA I]]"
$ws.Range("E17").Value = "HARD"
$ws.Range("F17").Value = "Implementação de projeto plant collections"
$ws.Range("G17").Value = 0.027777777777777776

# Match the selection left behind in the source workbook.
$ws.Range("F20").Select()
